$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new rows of data (row 3 and row 4), matching the existing
# rows' shape/content - same "Noun" method label, new elapsed/counts.
$ws.Cells.Item(3, 1).Value = 42600.829085648147
$ws.Cells.Item(3, 2).Value = "Noun"
$ws.Cells.Item(3, 3).Value = 6826
$ws.Cells.Item(3, 4).Value = 4017
$ws.Cells.Item(3, 5).Value = 622
$ws.Cells.Item(3, 6).Value = 103
$ws.Cells.Item(3, 7).Value = 54
$ws.Cells.Item(3, 8).Value = 64
$ws.Cells.Item(3, 9).Value = 34
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0

$ws.Cells.Item(4, 1).Value = 42600.881840277776
$ws.Cells.Item(4, 2).Value = "Noun"
$ws.Cells.Item(4, 3).Value = 8284
$ws.Cells.Item(4, 4).Value = 4017
$ws.Cells.Item(4, 5).Value = 622
$ws.Cells.Item(4, 6).Value = 103
$ws.Cells.Item(4, 7).Value = 54
$ws.Cells.Item(4, 8).Value = 64
$ws.Cells.Item(4, 9).Value = 34
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0

# Column A is formatted as a date/time, matching row 2's style
$ws.Range("A3:A4").NumberFormat = "m/d/yy h:mm"

# New (wider) date values make column A need a bit more room - widen it
# to fit, same as Excel's own best-fit behavior would after adding data.
$ws.Columns.Item(1).ColumnWidth = 14
